$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the notice date (cell is formatted with a date number format, so force
# text formatting first to keep the value stored as a literal text string)
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "2022-02-23"

# A15 and A16 were previously empty; now hold the recipient's name (on two lines)
$ws.Range("A15").Value = "Glee Star Enterprises"
$ws.Range("A16").Value = "Glee Star Enterprises`n"

# A17 previously held the address; now holds the new company's address
$ws.Range("A17").Value = "102 Centennial II Extension St., Pinagbuhatan, Pasig City`n"

# A20 salutation line
$ws.Range("A20").Value = "Dear Mr./Ms. Glee Star Enterprises"

# A23 main body paragraph referencing the contractor/company and the procurement purpose
$ws.Range("A23").Value = "Glee Star Enterprises`n that the Procurement of Supplies, Materials and Devices for the LOREM IPSUM shall commence upon receipt of the Notice to Proceed. "

# A43 signature line
$ws.Range("A43").Value = "                                    Glee Star Enterprises"
